$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 148, shifting existing rows 148-151 down to 149-152.
$ws.Rows("148:148").Insert()

# Fill in the values for the newly inserted row 148 (new weekly data point).
$ws.Cells.Item(148, 1).Value = 11
$ws.Cells.Item(148, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(148, 3).Value = "Bíobío"
$ws.Cells.Item(148, 4).Value = 45041
$ws.Cells.Item(148, 5).Value = 8
$ws.Cells.Item(148, 6).Value = 100112001
$ws.Cells.Item(148, 7).Value = "Berenjena"
$ws.Cells.Item(148, 8).Value = "Sin especificar"
$ws.Cells.Item(148, 9).Value = "Primera"
$ws.Cells.Item(148, 10).Value = 270
$ws.Cells.Item(148, 11).Value = 6500
$ws.Cells.Item(148, 12).Value = 7000
$ws.Cells.Item(148, 13).Value = 6778
$ws.Cells.Item(148, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(148, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(148, 16).Value = 113
$ws.Cells.Item(148, 17).Value = 60
$ws.Cells.Item(148, 18).Value = "Hortaliza"
